$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

$newRows = @(
    @(107, 1, "2024-06-17 05:15:49", 200, 11),
    @(108, 2, "2024-06-17 05:15:49", 200, 0)
)

$row = $lastRow + 1
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $row = $row + 1
}
